# Timesheet for Week 15 - date corrections
# Shift the weekly day labels (row 11-17, column A) forward by one week,
# and correct the "Week of:" date (G8) to match the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day-of-week labels for the new week (13/04 - 19/04)
$ws.Range("A11").Value = "Sun 13/04"
$ws.Range("A12").Value = "Mon 14/04"
$ws.Range("A13").Value = "Tue 15/04"
$ws.Range("A14").Value = "Wed 16/04"
$ws.Range("A15").Value = "Thur 17/04"
$ws.Range("A16").Value = "Fri   18/04"
$ws.Range("A17").Value = "Sat 19/04"

# "Week of:" date, corrected from 20/04/2014 to 13/04/2014
$ws.Range("G8").Value = 41742

# Restore the selection used when the sheet was last saved
$ws.Range("A17").Select()
